{"js": "// The author split the closing paragraph of the document (\"\u5b9a\u4e49\u76f8\u5173\u7684\u9519\u8bef\u7c7b\u578b\")\n// into two paragraphs: the original sentence stays on its own line, and a\n// brand-new second line (\"2 \u5173\u95ed\u9a71\u52a8\uff0c\u522b\u5fd8\u8bb0\u54af\") is appended after it, reusing\n// the original paragraph (so it keeps its paragraph mark / bookmark).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph right before the one we are duplicating/splitting,\n// i.e. the paragraph that holds \"\u4e0d\u662f\u8054\u901a\u548c\u5e7f\u4e1c\u79fb\u52a8\u7684\u8fd4\u56de \u63d0\u793a\u7c7b\u578b\uff1a\".\nconst items = paragraphs.items;\nlet anchor = null;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"\u4e0d\u662f\u8054\u901a\u548c\u5e7f\u4e1c\u79fb\u52a8\u7684\u8fd4\u56de\") !== -1) {\n    anchor = items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error(\"Anchor paragraph not found\");\n}\n\n// Insert a brand-new paragraph right after the anchor, carrying the text that\n// used to live in the final paragraph (\"\u5b9a\u4e49\u76f8\u5173\u7684\u9519\u8bef\u7c7b\u578b\"). Inserting it\n// \"after\" a plain paragraph (rather than \"before\" the target) keeps it free\n// of any inherited paragraph-level formatting, matching the original's\n// simple <w:p><w:r>... shape.\nanchor.insertParagraph(\"\u5b9a\u4e49\u76f8\u5173\u7684\u9519\u8bef\u7c7b\u578b\", Word.InsertLocation.after);\nawait context.sync();\n\n// Re-fetch paragraphs: the document's last paragraph is still the original\n// one (it now trails the newly inserted duplicate) and keeps its pPr /\n// bookmark. Replace its text with the new sentence.\nconst refreshed = body.paragraphs;\nrefreshed.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = refreshed.items[refreshed.items.length - 1];\nlastParagraph.insertText(\"2 \u5173\u95ed\u9a71\u52a8\uff0c\u522b\u5fd8\u8bb0\u54af\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The author split the closing paragraph of the document (\"\u5b9a\u4e49\u76f8\u5173\u7684\u9519\u8bef\u7c7b\u578b\")\n# into two paragraphs: the original sentence stays on its own line, and a\n# brand-new second line (\"2 \u5173\u95ed\u9a71\u52a8\uff0c\u522b\u5fd8\u8bb0\u54af\") is appended after it, reusing\n# the original paragraph (so it keeps its paragraph mark / bookmark \"_GoBack\").\n$d = $word.ActiveDocument\n\n# Locate the paragraph that will anchor the split: \"\u4e0d\u662f\u8054\u901a\u548c\u5e7f\u4e1c\u79fb\u52a8\u7684\u8fd4\u56de \u63d0\u793a\u7c7b\u578b\uff1a\".\n# (Find.Execute is used to confirm the text is present; the actual mutation is\n# then driven off $d.Paragraphs so the inserted paragraph lands exactly where\n# Word would put it, rather than via a Find range - which keeps growing to\n# swallow whatever is inserted at its own end.)\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"\u4e0d\u662f\u8054\u901a\u548c\u5e7f\u4e1c\u79fb\u52a8\u7684\u8fd4\u56de\")\nif (-not $found) {\n    throw \"Anchor paragraph '\u4e0d\u662f\u8054\u901a\u548c\u5e7f\u4e1c\u79fb\u52a8\u7684\u8fd4\u56de...' not found\"\n}\n\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.StartsWith(\"\u4e0d\u662f\u8054\u901a\u548c\u5e7f\u4e1c\u79fb\u52a8\u7684\u8fd4\u56de\")) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Anchor paragraph not found\"\n}\n\n# Insert a brand-new paragraph right after the anchor, carrying the text that\n# used to live in the final paragraph (\"\u5b9a\u4e49\u76f8\u5173\u7684\u9519\u8bef\u7c7b\u578b\"). Inserting via\n# InsertAfter on the anchor paragraph's own Range keeps the new paragraph\n# free of any inherited paragraph-level formatting, matching the original\n# simple <w:p><w:r>... shape.\n$anchorRange = $d.Paragraphs.Item($anchorIndex).Range\n$anchorRange.InsertAfter([char]13 + \"\u5b9a\u4e49\u76f8\u5173\u7684\u9519\u8bef\u7c7b\u578b\")\n\n# The document's last paragraph is still the original one (it now trails the\n# newly inserted duplicate) and keeps its pPr / bookmark. Replace its text,\n# excluding the trailing paragraph mark, with the new sentence.\n$lastParagraph = $d.Paragraphs.Last\n$lastRange = $lastParagraph.Range\n$lastRange.MoveEnd(1, -1) | Out-Null\n$lastRange.Text = \"2 \u5173\u95ed\u9a71\u52a8\uff0c\u522b\u5fd8\u8bb0\u54af\"\n"}
